$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert header row at the top; existing data (rows 2-6) stays where it is.
$ws.Range("B1").Value = "WORDS"
$ws.Range("C1").Value = "TYPE"
$ws.Range("D1").Value = "TR"

# Move the selection to match the target state.
$ws.Range("E5").Select()
